$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# Handle the value that collides first: 94÷6= -> 85÷2=
Replace-Text "94÷6=" "85÷2="

Replace-Text "33÷7=" "25÷3="
Replace-Text "70÷3=" "63÷6="
Replace-Text "73÷5=" "64÷5="
Replace-Text "52÷9=" "16÷9="
Replace-Text "61÷3=" "87÷4="
Replace-Text "79÷7=" "55÷5="
Replace-Text "23÷7=" "92÷2="
Replace-Text "42÷5=" "97÷7="
Replace-Text "30÷2=" "64÷9="
Replace-Text "78÷6=" "41÷3="
Replace-Text "27÷7=" "94÷6="
Replace-Text "57÷6=" "91÷2="
Replace-Text "19÷9=" "34÷4="
Replace-Text "50÷4=" "89÷9="
Replace-Text "14÷2=" "27÷9="
Replace-Text "66÷5=" "65÷8="
Replace-Text "50÷8=" "71÷9="
Replace-Text "86÷3=" "12÷6="
Replace-Text "69÷3=" "70÷2="
Replace-Text "43÷7=" "93÷5="
Replace-Text "81÷5=" "74÷5="
Replace-Text "41÷2=" "48÷3="
Replace-Text "30÷7=" "38÷4="
Replace-Text "14÷6=" "35÷9="
